$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[0.26096798514997205, 0.3984326536276695]"
$ws.Range("M2").Value = 0.000000002021911749139349
$ws.Range("N2").Value = 0.000000002021911749139349
$ws.Range("P2").Value = "[-1.2578949563923096, -0.8553685703467693]"
$ws.Range("T2").Value = "[0.2759859091735271, 0.34766814850615374]"
$ws.Range("U2").Value = 0.00000000000002753353101070388
$ws.Range("V2").Value = 0.00000000000002753353101070388
$ws.Range("X2").Value = 3.257737737737815
$ws.Range("Y2").Value = 4.790790790790912
